# Updated symbol list on Tue Feb  7 21:46:15 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) figures for the crypto ticker rows.
# The source cells are plain text (t="inlineStr") — e.g. "331.77", "0.84%" —
# so each new value is written with a leading apostrophe to force Excel to
# store it as a text literal (quote-prefix) instead of auto-coercing the
# numeric-looking strings into Number cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.89"
$ws.Range("E2").Value = "'0.87%"

$ws.Range("D3").Value = "'44.68"
$ws.Range("E3").Value = "'1.39%"

$ws.Range("D4").Value = "'5.550"

$ws.Range("D5").Value = "'0.08210"
$ws.Range("E5").Value = "'1.79%"

$ws.Range("D6").Value = "'2.054"
$ws.Range("E6").Value = "'2.10%"

$ws.Range("D7").Value = "'0.9773"
$ws.Range("E7").Value = "'2.61%"

$ws.Range("D8").Value = "'0.1120"
$ws.Range("E8").Value = "'-3.23%"

$ws.Range("D9").Value = "'0.1909"
$ws.Range("E9").Value = "'3.05%"

$ws.Range("D10").Value = "'10.24"
$ws.Range("E10").Value = "'-13.68%"

$ws.Range("E11").Value = "'1.42%"

$ws.Range("D12").Value = "'0.04707"
$ws.Range("E12").Value = "'-0.97%"

$ws.Range("E13").Value = "'-0.94%"

$ws.Range("D14").Value = "'0.001259"
$ws.Range("E14").Value = "'-2.08%"

$ws.Range("D15").Value = "'0.04109"
$ws.Range("E15").Value = "'-2.96%"

$ws.Range("D16").Value = "'0.005919"
$ws.Range("E16").Value = "'-0.03%"

$ws.Range("D17").Value = "'3.350"
$ws.Range("E17").Value = "'-0.62%"

$ws.Range("D18").Value = "'4.429"
$ws.Range("E18").Value = "'2.30%"

$ws.Range("E19").Value = "'2.79%"

$ws.Range("E20").Value = "'-3.54%"

$ws.Range("D21").Value = "'0.1376"
$ws.Range("E21").Value = "'-2.38%"

$ws.Range("D22").Value = "'0.2491"
$ws.Range("E22").Value = "'-0.70%"

$ws.Range("D23").Value = "'0.001301"
$ws.Range("E23").Value = "'3.96%"

$ws.Range("D24").Value = "'0.004390"
$ws.Range("E24").Value = "'1.73%"

$ws.Range("D25").Value = "'0.0001279"
$ws.Range("E25").Value = "'7.31%"

$ws.Range("D26").Value = "'0.0003742"
$ws.Range("E26").Value = "'-5.96%"

$ws.Range("D38").Value = "'0.02766"
$ws.Range("E38").Value = "'4.98%"

$ws.Range("D39").Value = "'0.05727"
$ws.Range("E39").Value = "'3.23%"

$ws.Range("D40").Value = "'0.007640"
$ws.Range("E40").Value = "'0.95%"

$ws.Range("D41").Value = "'0.1425"
$ws.Range("E41").Value = "'1.34%"

$ws.Range("D42").Value = "'0.007533"
$ws.Range("E42").Value = "'-6.88%"

$ws.Range("D43").Value = "'0.001972"
$ws.Range("E43").Value = "'-2.29%"

$ws.Range("D44").Value = "'0.008294"
$ws.Range("E44").Value = "'-1.03%"

$ws.Range("D45").Value = "'0.00007029"
$ws.Range("E45").Value = "'-4.08%"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.09%"

$ws.Range("D47").Value = "'0.0005802"
$ws.Range("E47").Value = "'-0.16%"

$ws.Range("D48").Value = "'0.003574"
$ws.Range("E48").Value = "'-26.43%"

$ws.Range("D49").Value = "'0.002521"
$ws.Range("E49").Value = "'9.53%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.09%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.09%"
